$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a brand-new paragraph right after the "settings button ...
#    calendar activity" paragraph (paragraph 5) and before the
#    "database.h" paragraph, describing the three new methods added to
#    the main page.
# ---------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*calendar activity as soon as possible.*") {
        $anchor = $p
        break
    }
}

# Create a fresh, empty paragraph right after the anchor paragraph.
$anchor.Range.InsertParagraphAfter()

# Re-locate that freshly minted paragraph (it now directly follows the
# anchor paragraph) so we can fill it in.
$newParaIndex = $anchor.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)

$newText = "The main page has three new methods: onCreateOptionsMenu, onOptionsItemSelected, and onClick. The onCreateOptionsMenu adds a button to the top right for adding a task. The onOptionsItemSelected method adds the " + [char]0x201C + "add or " + [char]0x201C + "cancel" + [char]0x201D + " option after a task has been input. The onClick method adds a settings button in the shape of a gear cog. This method takes you to a separate page with a list of options for app settings. There is a new field called mTaskList. It is a ListView object and is used for displaying the task list ArrayList. The add a task functionality still needs additional task fields added to the create task menu screen. A UI update method needs to be added to refresh the lask list screen when a new task is added. Along with the settings button a history button needs to be added to take the user to a screen to view a list of completed or canceled tasks. It will display similar to the task list ArrayList and ListView."

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xmlFragment)

# ---------------------------------------------------------------------
# 2) Resize the inline diagram picture.
# ---------------------------------------------------------------------

$shp = $d.InlineShapes.Item(1)
$shp.Width = 6742190 / 914400 * 72
$shp.Height = 3532584 / 914400 * 72
